# DoE_TempFlume.xlsx edit:
#  - "Tests Scale 1-21" sheet: drop the unused trailing columns (G:N) that held
#    "Unnamed: 7..12" headers and stray turbulence-model scratch values, and
#    clear the "Crashed" marker out of rows 3-38 (only row 2 keeps it, matching
#    the target workbook).
#  - "Tests Scale 1-25" sheet: populate the "Crashed"/"Completed" results
#    column (E) for every test row, monitoring which runs finished vs which
#    ones crashed.

$wb = $excel.ActiveWorkbook

# --- "Tests Scale 1-21" -----------------------------------------------------
$ws21 = $wb.Worksheets.Item("Tests Scale 1-21")

# Remove the stray/unused columns G through N entirely (not just their
# values) so the sheet's used range shrinks back down to A:F.
$ws21.Range("G1:N1").EntireColumn.Delete()

# Only the very first data row (row 2) keeps its "Crashed" result; the rest
# of the legacy "Crashed" markers in rows 3-38 are cleared out.
$ws21.Range("E3:E38").ClearContents()

# --- "Tests Scale 1-25" -----------------------------------------------------
$ws25 = $wb.Worksheets.Item("Tests Scale 1-25")

# Mark every test row as "Crashed" by default...
$ws25.Range("E2:E61").Value = "Crashed"

# ...except the couple of runs that actually finished, which get "Completed".
$ws25.Range("E22").Value = "Completed"
$ws25.Range("E30").Value = "Completed"
